# "#5: cash & deposit done"
# Extend the deposit (存款) sheet so that it carries the same
# property/category/legislator metadata columns already present on the
# land (土地) and stock (股票) sheets, and turn its first row into a
# proper header row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 土地 (land)   - used as a safe source for the date string
$ws2 = $wb.Worksheets.Item(2)   # 存款 (deposit) - sheet being edited

$lastRow = 13

# ------------------------------------------------------------------
# 1. Fix up the header row (row 1).
#    E1 used to hold the data value "謝國樑"; it should hold the
#    column label "owner" instead. F1 used to hold a bare numeric
#    literal (1); it should hold the column label "total" instead.
# ------------------------------------------------------------------
$ws2.Range("E1").Value = "owner"
$ws2.Range("F1").Value = "total"

# New header labels for the newly added columns G1:M1, matching the
# layout already used on the other two sheets.
$ws2.Range("G1").Value = "property_category"
$ws2.Range("H1").Value = "category"
$ws2.Range("I1").Value = "date"
$ws2.Range("J1").Value = "legislator_name"
$ws2.Range("K1").Value = "legislator_id"
$ws2.Range("L1").Value = "source_file"
$ws2.Range("M1").Value = "index"

# ------------------------------------------------------------------
# 2. Populate the new data columns (G:M) for every data row.
#    G = property_category -> "deposit"
#    H = category           -> "normal"
#    I = date                -> "2012-05-01" (copied from sheet1 so the
#                                text isn't reinterpreted as a real date)
#    J = legislator_name    -> "謝國樑"
#    K = legislator_id      -> 1387
#    L = source_file        -> "tmpa28e1"
#    M = index               -> same value as column A on that row
# ------------------------------------------------------------------

# Copy the already-existing date string cell from sheet1 (row 2, column K)
# so that Excel keeps it as literal text instead of converting it to a
# serial date number.
$ws1.Range("K2").Copy() | Out-Null

for ($r = 2; $r -le $lastRow; $r++) {
    $ws2.Range("G$r").Value = "deposit"
    $ws2.Range("H$r").Value = "normal"
    $ws2.Range("I$r").PasteSpecial(-4163) | Out-Null   # xlPasteAll (keeps text "2012-05-01")
    $ws2.Range("J$r").Value = "謝國樑"
    $ws2.Range("K$r").Value = 1387
    $ws2.Range("L$r").Value = "tmpa28e1"
    $ws2.Range("M$r").Value = $ws2.Range("A$r").Value2
}
